$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.077.60'
$ws.Range("E2").Value = '  +1.57%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.888.11'
$ws.Range("E3").Value = '  +1.05%  '
$ws.Range("E4").Value = '  +1.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '336.15'
$ws.Range("E5").Value = '  +1.50%  '
$ws.Range("E6").Value = '  +0.94%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4762'
$ws.Range("E7").Value = '  +1.57%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3955'
$ws.Range("E8").Value = '  +0.52%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.12'
$ws.Range("E9").Value = '  -1.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08039'
$ws.Range("E10").Value = '  -0.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.021'
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.92'
$ws.Range("E12").Value = '  +0.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.907.69'
$ws.Range("E13").Value = '  +3.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.050'
$ws.Range("E14").Value = '  +1.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.216'
$ws.Range("E15").Value = '  +1.12%  '
$ws.Range("E16").Value = '  +0.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.61'
$ws.Range("E17").Value = '  +2.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06760'
$ws.Range("E18").Value = '  +2.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001051'
$ws.Range("E19").Value = '  +0.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.09'
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.011'
$ws.Range("E21").Value = '  +0.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '28.051.67'
$ws.Range("E22").Value = '  +1.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.522'
$ws.Range("E23").Value = '  +0.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.04'
$ws.Range("E24").Value = '  +0.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.353'
$ws.Range("E25").Value = '  +1.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.124.52'
$ws.Range("E26").Value = '  +2.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.73'
$ws.Range("E27").Value = '  +0.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.03'
$ws.Range("E28").Value = '  -0.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.108'
$ws.Range("E29").Value = '  +0.89%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.521'
$ws.Range("E30").Value = '  -0.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '121.86'
$ws.Range("E31").Value = '  -0.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9782'
$ws.Range("E32").Value = '  +1.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09597'
$ws.Range("E33").Value = '  +1.19%  '
$ws.Range("E34").Value = '  +1.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.364'
$ws.Range("E35").Value = '  +0.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.366'
$ws.Range("E36").Value = '  -5.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02255'
$ws.Range("E37").Value = '  +0.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06076'
$ws.Range("E38").Value = '  -0.16%  '
$ws.Range("E39").Value = '  -1.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.219'
$ws.Range("E40").Value = '  +1.37%  '
$ws.Range("E41").Value = '  +0.79%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5985'
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1896'
$ws.Range("E43").Value = '  +0.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.36'
$ws.Range("E44").Value = '  +1.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.263'
$ws.Range("E45").Value = '  +0.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5667'
$ws.Range("E46").Value = '  -0.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.24'
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.936'
$ws.Range("E48").Value = '  +0.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.360'
$ws.Range("E49").Value = '  -0.79%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06827'
$ws.Range("E50").Value = '  -0.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '112.38'
$ws.Range("E51").Value = '  -1.64%  '
